# Update the date heading and every arithmetic-problem answer cell in the
# table to reflect the regenerated "output at 596fc94" values.
$d = $word.ActiveDocument

$d.Content.Find.Execute("2023-06-15 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-06-16 Friday", 2) | Out-Null
$d.Content.Find.Execute("76-64=12", $true, $false, $false, $false, $false, $true, 1, $false, "71-63=8", 2) | Out-Null
$d.Content.Find.Execute("54+9=63", $true, $false, $false, $false, $false, $true, 1, $false, "57+5=62", 2) | Out-Null
$d.Content.Find.Execute("95-24=71", $true, $false, $false, $false, $false, $true, 1, $false, "91-38=53", 2) | Out-Null
$d.Content.Find.Execute("39+53=92", $true, $false, $false, $false, $false, $true, 1, $false, "28+7=35", 2) | Out-Null
$d.Content.Find.Execute("91-60=31", $true, $false, $false, $false, $false, $true, 1, $false, "39-20=19", 2) | Out-Null
$d.Content.Find.Execute("69-5=64", $true, $false, $false, $false, $false, $true, 1, $false, "73+7=80", 2) | Out-Null
$d.Content.Find.Execute("53-3=50", $true, $false, $false, $false, $false, $true, 1, $false, "64-53=11", 2) | Out-Null
$d.Content.Find.Execute("41+40=81", $true, $false, $false, $false, $false, $true, 1, $false, "48-1=47", 2) | Out-Null
$d.Content.Find.Execute("20+11=31", $true, $false, $false, $false, $false, $true, 1, $false, "83+9=92", 2) | Out-Null
$d.Content.Find.Execute("10+79=89", $true, $false, $false, $false, $false, $true, 1, $false, "73-56=17", 2) | Out-Null
$d.Content.Find.Execute("24+40=64", $true, $false, $false, $false, $false, $true, 1, $false, "6+82=88", 2) | Out-Null
$d.Content.Find.Execute("78-12=66", $true, $false, $false, $false, $false, $true, 1, $false, "38-36=2", 2) | Out-Null
$d.Content.Find.Execute("23+37=60", $true, $false, $false, $false, $false, $true, 1, $false, "58-28=30", 2) | Out-Null
$d.Content.Find.Execute("72-60=12", $true, $false, $false, $false, $false, $true, 1, $false, "44+45=89", 2) | Out-Null
$d.Content.Find.Execute("59-56=3", $true, $false, $false, $false, $false, $true, 1, $false, "18+57=75", 2) | Out-Null
$d.Content.Find.Execute("54-45=9", $true, $false, $false, $false, $false, $true, 1, $false, "24+31=55", 2) | Out-Null
$d.Content.Find.Execute("70+6=76", $true, $false, $false, $false, $false, $true, 1, $false, "53+1=54", 2) | Out-Null
$d.Content.Find.Execute("48-14=34", $true, $false, $false, $false, $false, $true, 1, $false, "69-18=51", 2) | Out-Null
$d.Content.Find.Execute("52+21=73", $true, $false, $false, $false, $false, $true, 1, $false, "55+0=55", 2) | Out-Null
$d.Content.Find.Execute("49+1=50", $true, $false, $false, $false, $false, $true, 1, $false, "4+8=12", 2) | Out-Null
$d.Content.Find.Execute("26-1=25", $true, $false, $false, $false, $false, $true, 1, $false, "45+47=92", 2) | Out-Null
$d.Content.Find.Execute("89-51=38", $true, $false, $false, $false, $false, $true, 1, $false, "85-23=62", 2) | Out-Null
$d.Content.Find.Execute("98-12=86", $true, $false, $false, $false, $false, $true, 1, $false, "89-42=47", 2) | Out-Null
$d.Content.Find.Execute("60+16=76", $true, $false, $false, $false, $false, $true, 1, $false, "61-36=25", 2) | Out-Null
$d.Content.Find.Execute("4+71=75", $true, $false, $false, $false, $false, $true, 1, $false, "42-22=20", 2) | Out-Null
$d.Content.Find.Execute("63+34=97", $true, $false, $false, $false, $false, $true, 1, $false, "80-19=61", 2) | Out-Null
$d.Content.Find.Execute("26+16=42", $true, $false, $false, $false, $false, $true, 1, $false, "18+40=58", 2) | Out-Null
$d.Content.Find.Execute("56-35=21", $true, $false, $false, $false, $false, $true, 1, $false, "69+6=75", 2) | Out-Null
$d.Content.Find.Execute("55-19=36", $true, $false, $false, $false, $false, $true, 1, $false, "62+25=87", 2) | Out-Null
$d.Content.Find.Execute("16+4=20", $true, $false, $false, $false, $false, $true, 1, $false, "31+2=33", 2) | Out-Null
$d.Content.Find.Execute("11+39=50", $true, $false, $false, $false, $false, $true, 1, $false, "52-35=17", 2) | Out-Null
$d.Content.Find.Execute("5+57=62", $true, $false, $false, $false, $false, $true, 1, $false, "67-33=34", 2) | Out-Null
$d.Content.Find.Execute("90-28=62", $true, $false, $false, $false, $false, $true, 1, $false, "46-8=38", 2) | Out-Null
$d.Content.Find.Execute("76-13=63", $true, $false, $false, $false, $false, $true, 1, $false, "80-37=43", 2) | Out-Null
$d.Content.Find.Execute("53-29=24", $true, $false, $false, $false, $false, $true, 1, $false, "5+24=29", 2) | Out-Null
$d.Content.Find.Execute("74-66=8", $true, $false, $false, $false, $false, $true, 1, $false, "3+15=18", 2) | Out-Null
$d.Content.Find.Execute("42+47=89", $true, $false, $false, $false, $false, $true, 1, $false, "9+31=40", 2) | Out-Null
$d.Content.Find.Execute("12+4=16", $true, $false, $false, $false, $false, $true, 1, $false, "40+9=49", 2) | Out-Null
$d.Content.Find.Execute("76+5=81", $true, $false, $false, $false, $false, $true, 1, $false, "1+67=68", 2) | Out-Null
$d.Content.Find.Execute("98-16=82", $true, $false, $false, $false, $false, $true, 1, $false, "70+26=96", 2) | Out-Null
$d.Content.Find.Execute("31+47=78", $true, $false, $false, $false, $false, $true, 1, $false, "41-9=32", 2) | Out-Null
$d.Content.Find.Execute("42-41=1", $true, $false, $false, $false, $false, $true, 1, $false, "90+9=99", 2) | Out-Null
$d.Content.Find.Execute("89-50=39", $true, $false, $false, $false, $false, $true, 1, $false, "2+87=89", 2) | Out-Null
$d.Content.Find.Execute("68-60=8", $true, $false, $false, $false, $false, $true, 1, $false, "94-41=53", 2) | Out-Null
$d.Content.Find.Execute("94-62=32", $true, $false, $false, $false, $false, $true, 1, $false, "98-54=44", 2) | Out-Null
$d.Content.Find.Execute("14+16=30", $true, $false, $false, $false, $false, $true, 1, $false, "30+63=93", 2) | Out-Null
$d.Content.Find.Execute("12+58=70", $true, $false, $false, $false, $false, $true, 1, $false, "10+43=53", 2) | Out-Null
$d.Content.Find.Execute("40+56=96", $true, $false, $false, $false, $false, $true, 1, $false, "82-23=59", 2) | Out-Null
$d.Content.Find.Execute("12+42=54", $true, $false, $false, $false, $false, $true, 1, $false, "10-2=8", 2) | Out-Null
$d.Content.Find.Execute("89-83=6", $true, $false, $false, $false, $false, $true, 1, $false, "16+76=92", 2) | Out-Null
$d.Content.Find.Execute("47+29=76", $true, $false, $false, $false, $false, $true, 1, $false, "53-45=8", 2) | Out-Null
$d.Content.Find.Execute("76-66=10", $true, $false, $false, $false, $false, $true, 1, $false, "20+72=92", 2) | Out-Null
$d.Content.Find.Execute("16-15=1", $true, $false, $false, $false, $false, $true, 1, $false, "33+55=88", 2) | Out-Null
$d.Content.Find.Execute("5+50=55", $true, $false, $false, $false, $false, $true, 1, $false, "37+20=57", 2) | Out-Null
$d.Content.Find.Execute("69-65=4", $true, $false, $false, $false, $false, $true, 1, $false, "79-7=72", 2) | Out-Null
$d.Content.Find.Execute("55-46=9", $true, $false, $false, $false, $false, $true, 1, $false, "24+60=84", 2) | Out-Null
$d.Content.Find.Execute("44+43=87", $true, $false, $false, $false, $false, $true, 1, $false, "91-9=82", 2) | Out-Null
$d.Content.Find.Execute("99-34=65", $true, $false, $false, $false, $false, $true, 1, $false, "41+14=55", 2) | Out-Null
$d.Content.Find.Execute("29+20=49", $true, $false, $false, $false, $false, $true, 1, $false, "4+78=82", 2) | Out-Null
$d.Content.Find.Execute("58+16=74", $true, $false, $false, $false, $false, $true, 1, $false, "51+24=75", 2) | Out-Null
$d.Content.Find.Execute("92-43=49", $true, $false, $false, $false, $false, $true, 1, $false, "71-45=26", 2) | Out-Null
$d.Content.Find.Execute("79-31=48", $true, $false, $false, $false, $false, $true, 1, $false, "50+17=67", 2) | Out-Null
$d.Content.Find.Execute("8+85=93", $true, $false, $false, $false, $false, $true, 1, $false, "7+56=63", 2) | Out-Null
$d.Content.Find.Execute("10+65=75", $true, $false, $false, $false, $false, $true, 1, $false, "20+2=22", 2) | Out-Null
$d.Content.Find.Execute("50-31=19", $true, $false, $false, $false, $false, $true, 1, $false, "4+20=24", 2) | Out-Null
$d.Content.Find.Execute("6+74=80", $true, $false, $false, $false, $false, $true, 1, $false, "75-63=12", 2) | Out-Null
$d.Content.Find.Execute("39+12=51", $true, $false, $false, $false, $false, $true, 1, $false, "1+44=45", 2) | Out-Null
$d.Content.Find.Execute("29+61=90", $true, $false, $false, $false, $false, $true, 1, $false, "37-23=14", 2) | Out-Null
$d.Content.Find.Execute("35-14=21", $true, $false, $false, $false, $false, $true, 1, $false, "26+50=76", 2) | Out-Null
$d.Content.Find.Execute("47-35=12", $true, $false, $false, $false, $false, $true, 1, $false, "60-55=5", 2) | Out-Null
$d.Content.Find.Execute("5+75=80", $true, $false, $false, $false, $false, $true, 1, $false, "74+16=90", 2) | Out-Null
$d.Content.Find.Execute("78-75=3", $true, $false, $false, $false, $false, $true, 1, $false, "99-11=88", 2) | Out-Null
$d.Content.Find.Execute("15+8=23", $true, $false, $false, $false, $false, $true, 1, $false, "49-18=31", 2) | Out-Null
$d.Content.Find.Execute("91-78=13", $true, $false, $false, $false, $false, $true, 1, $false, "81+15=96", 2) | Out-Null
$d.Content.Find.Execute("13+14=27", $true, $false, $false, $false, $false, $true, 1, $false, "47+7=54", 2) | Out-Null
$d.Content.Find.Execute("11+53=64", $true, $false, $false, $false, $false, $true, 1, $false, "26-14=12", 2) | Out-Null
$d.Content.Find.Execute("3+28=31", $true, $false, $false, $false, $false, $true, 1, $false, "31-17=14", 2) | Out-Null
$d.Content.Find.Execute("89-89=0", $true, $false, $false, $false, $false, $true, 1, $false, "2+41=43", 2) | Out-Null
$d.Content.Find.Execute("72-61=11", $true, $false, $false, $false, $false, $true, 1, $false, "7+23=30", 2) | Out-Null
$d.Content.Find.Execute("9+42=51", $true, $false, $false, $false, $false, $true, 1, $false, "20+79=99", 2) | Out-Null
$d.Content.Find.Execute("35+37=72", $true, $false, $false, $false, $false, $true, 1, $false, "35+38=73", 2) | Out-Null
$d.Content.Find.Execute("48+41=89", $true, $false, $false, $false, $false, $true, 1, $false, "29+34=63", 2) | Out-Null
$d.Content.Find.Execute("32+21=53", $true, $false, $false, $false, $false, $true, 1, $false, "74-10=64", 2) | Out-Null
$d.Content.Find.Execute("71-19=52", $true, $false, $false, $false, $false, $true, 1, $false, "60-49=11", 2) | Out-Null
$d.Content.Find.Execute("28+27=55", $true, $false, $false, $false, $false, $true, 1, $false, "37-34=3", 2) | Out-Null
$d.Content.Find.Execute("92-58=34", $true, $false, $false, $false, $false, $true, 1, $false, "58-6=52", 2) | Out-Null
$d.Content.Find.Execute("48-30=18", $true, $false, $false, $false, $false, $true, 1, $false, "78-0=78", 2) | Out-Null
$d.Content.Find.Execute("30+68=98", $true, $false, $false, $false, $false, $true, 1, $false, "75-29=46", 2) | Out-Null
$d.Content.Find.Execute("96-6=90", $true, $false, $false, $false, $false, $true, 1, $false, "64-17=47", 2) | Out-Null
$d.Content.Find.Execute("32+17=49", $true, $false, $false, $false, $false, $true, 1, $false, "96-73=23", 2) | Out-Null
$d.Content.Find.Execute("81+2=83", $true, $false, $false, $false, $false, $true, 1, $false, "94-55=39", 2) | Out-Null
$d.Content.Find.Execute("69-3=66", $true, $false, $false, $false, $false, $true, 1, $false, "3+22=25", 2) | Out-Null
$d.Content.Find.Execute("4+45=49", $true, $false, $false, $false, $false, $true, 1, $false, "6+56=62", 2) | Out-Null
$d.Content.Find.Execute("49+4=53", $true, $false, $false, $false, $false, $true, 1, $false, "3+22=25", 2) | Out-Null
$d.Content.Find.Execute("79+5=84", $true, $false, $false, $false, $false, $true, 1, $false, "2+34=36", 2) | Out-Null
$d.Content.Find.Execute("49+32=81", $true, $false, $false, $false, $false, $true, 1, $false, "5+60=65", 2) | Out-Null
$d.Content.Find.Execute("15+28=43", $true, $false, $false, $false, $false, $true, 1, $false, "60-43=17", 2) | Out-Null
$d.Content.Find.Execute("22+18=40", $true, $false, $false, $false, $false, $true, 1, $false, "57-37=20", 2) | Out-Null
$d.Content.Find.Execute("19-7=12", $true, $false, $false, $false, $false, $true, 1, $false, "78-59=19", 2) | Out-Null
$d.Content.Find.Execute("58-37=21", $true, $false, $false, $false, $false, $true, 1, $false, "60+13=73", 2) | Out-Null
